$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44326, 1, 5, 109.051254089422),
    @(44327, 0, 5, 109.051254089422),
    @(44328, 0, 5, 109.051254089422),
    @(44329, 2, 4, 87.24100327153762)
)

$row = 252
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]

    $ws.Cells.Item(251, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $row++
}
